$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "contratante"
$ws.Range("C2:C16").Value = "KFC"

$ws.Range("C2:C16").Select()
